$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column I
$ws.Range("I1").Value = "Số giờ tăng ca"

# Fill column I (rows 2-21) with 0, matching the other numeric columns
for ($r = 2; $r -le 21; $r++) {
    $ws.Cells.Item($r, 9).Value = 0
}

# Update the selected cell to reflect the new active cell/range in the edited file
$ws.Range("M5").Select()
